$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Config Name and data values in row 2
$ws.Range("A2").Value = "Test"
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

# F2 keeps its original "quote-prefixed number" formatting (style index 1)
# even though its value changed; re-apply that formatting from F3, which
# carries the same style untouched.
$ws.Range("F3").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to A2
[void]$ws.Range("A2").Select()
